# Removed Test Case Inter-Dependency
# - Decouple the "shortname" value from the numeric product code (2560) by
#   replacing it with a literal, non-dependent short code ("256d").
# - Rename the product name suffix from "-UPFRONT" to "-UPF-1st".
# - Reset the saved selection on the input sheet back to the top (B1) instead
#   of wherever the last test run left it (B34), and make the output sheet
#   the active tab/selection instead of the input sheet.

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update the product name text (B1) on both sheets (the output sheet mirrors
# the same product-name string so they continue to share the pooled string).
$wsInput.Range("B1").Value  = "2560-MS-EI-DB-SAR-REC-NOCOM-RNI-CTPD-SAR-MD-TR-2-DATE-VAR-INST-UPF-1st"
$wsOutput.Range("B1").Value = "2560-MS-EI-DB-SAR-REC-NOCOM-RNI-CTPD-SAR-MD-TR-2-DATE-VAR-INST-UPF-1st"

# Update the short name (B2) on the input sheet: was the bare number 2560,
# now a literal string so it no longer depends on another test case's state.
$wsInput.Range("B2").Value = "256d"

# Reset the remembered selection on the input sheet to B1 (was B34).
$wsInput.Range("B1").Select()

# Make the output sheet the active tab/selection, deselecting the input tab.
$wsOutput.Activate()
$wsOutput.Range("B1").Select()
